$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organizations")

# --- Address (Street_Address, column E) corrections: strip city/state/zip, keep street only ---
$ws.Range('E2').Value = '18000 Pacific Coast Hwy'
$ws.Range('E3').Value = '587 Palm Canyon Dr #182'
$ws.Range('E4').Value = '1916 Albans Road'
$ws.Range('E9').Value = '1916 Albans Road'
$ws.Range('E10').Value = '1916 Albans Road'
$ws.Range('E13').Value = '1916 Albans Road'
$ws.Range('E14').Value = '1916 Albans Road'
$ws.Range('E15').Value = '1916 Albans Road'
$ws.Range('E16').Value = '1916 Albans Road'
$ws.Range('E17').Value = '1916 Albans Road'
$ws.Range('E18').Value = '1916 Albans Road'
$ws.Range('E21').Value = '1916 Albans Road'
$ws.Range('E22').Value = '1916 Albans Road'
$ws.Range('E23').Value = '1916 Albans Road'
$ws.Range('E25').Value = '1916 Albans Road'
$ws.Range('E28').Value = '1916 Albans Road'
$ws.Range('E29').Value = '1916 Albans Road'
$ws.Range('E31').Value = '1916 Albans Road'
$ws.Range('E36').Value = '1916 Albans Road'
$ws.Range('E37').Value = '1916 Albans Road'
$ws.Range('E40').Value = '1916 Albans Road'
$ws.Range('E42').Value = '1916 Albans Road'
$ws.Range('E44').Value = '1916 Albans Road'
$ws.Range('E45').Value = '1916 Albans Road'
$ws.Range('E47').Value = '1916 Albans Road'
$ws.Range('E51').Value = '1916 Albans Road'
$ws.Range('E55').Value = '1916 Albans Road'
$ws.Range('E57').Value = '1916 Albans Road'
$ws.Range('E58').Value = '1916 Albans Road'
$ws.Range('E61').Value = '1916 Albans Road'
$ws.Range('E63').Value = '1916 Albans Road'
$ws.Range('E64').Value = '1916 Albans Road'
$ws.Range('E67').Value = '1916 Albans Road'
$ws.Range('E70').Value = '1916 Albans Road'
$ws.Range('E72').Value = '1916 Albans Road'
$ws.Range('E73').Value = '1916 Albans Road'
$ws.Range('E74').Value = '1916 Albans Road'
$ws.Range('E80').Value = '1916 Albans Road'
$ws.Range('E81').Value = '1916 Albans Road'
$ws.Range('E85').Value = '1916 Albans Road'
$ws.Range('E86').Value = '1916 Albans Road'
$ws.Range('E87').Value = '1916 Albans Road'
$ws.Range('E88').Value = '1916 Albans Road'
$ws.Range('E92').Value = '1916 Albans Road'
$ws.Range('E93').Value = '1916 Albans Road'
$ws.Range('E95').Value = '1916 Albans Road'
$ws.Range('E98').Value = '1916 Albans Road'
$ws.Range('E99').Value = '1916 Albans Road'
$ws.Range('E101').Value = '1916 Albans Road'
$ws.Range('E105').Value = '1916 Albans Road'
$ws.Range('E5').Value = '4900 Shoreline Hwy '
$ws.Range('E6').Value = '1250 Addison St #101 '
$ws.Range('E7').Value = '3521 14 Mile House Rd '
$ws.Range('E8').Value = '3842 Warner Ave '
$ws.Range('E11').Value = '715 P Street'
$ws.Range('E12').Value = '6605 San Leandro St '
$ws.Range('E19').Value = '2950 Peralta Oaks Ct'
$ws.Range('E20').Value = '2810 Pio Pico Dr'
$ws.Range('E24').Value = '101 Cooper St '
$ws.Range('E26').Value = '3140 Sierrama Dr '
$ws.Range('E27').Value = 'Post Office Box 829 '
$ws.Range('E30').Value = '51500 CA-74 '
$ws.Range('E32').Value = '604 Main St '

# --- Date-like cells (columns AE:AJ) converted from date serials to text "1-19%" ---
$ws.Range('AH2').Value = '1-19%'
$ws.Range('AJ2').Value = '1-19%'
$ws.Range('AE3').Value = '1-19%'
$ws.Range('AF3').Value = '1-19%'
$ws.Range('AG3').Value = '1-19%'
$ws.Range('AF7').Value = '1-19%'
$ws.Range('AE8').Value = '1-19%'
$ws.Range('AG8').Value = '1-19%'
$ws.Range('AH8').Value = '1-19%'
$ws.Range('AE14').Value = '1-19%'
$ws.Range('AH14').Value = '1-19%'
$ws.Range('AF15').Value = '1-19%'
$ws.Range('AG18').Value = '1-19%'
$ws.Range('AG20').Value = '1-19%'
$ws.Range('AE21').Value = '1-19%'
$ws.Range('AG22').Value = '1-19%'
$ws.Range('AH22').Value = '1-19%'
$ws.Range('AF23').Value = '1-19%'
$ws.Range('AH23').Value = '1-19%'
$ws.Range('AG24').Value = '1-19%'
$ws.Range('AF30').Value = '1-19%'
$ws.Range('AG30').Value = '1-19%'
$ws.Range('AJ30').Value = '1-19%'
$ws.Range('AG32').Value = '1-19%'
$ws.Range('AJ32').Value = '1-19%'
$ws.Range('AJ33').Value = '1-19%'
$ws.Range('AG34').Value = '1-19%'
$ws.Range('AE36').Value = '1-19%'
$ws.Range('AG36').Value = '1-19%'
$ws.Range('AH36').Value = '1-19%'
$ws.Range('AH37').Value = '1-19%'
$ws.Range('AI37').Value = '1-19%'
$ws.Range('AJ39').Value = '1-19%'
$ws.Range('AI40').Value = '1-19%'
$ws.Range('AH41').Value = '1-19%'
$ws.Range('AI41').Value = '1-19%'
$ws.Range('AE43').Value = '1-19%'
$ws.Range('AF43').Value = '1-19%'
$ws.Range('AE44').Value = '1-19%'
$ws.Range('AE45').Value = '1-19%'
$ws.Range('AG45').Value = '1-19%'
$ws.Range('AE47').Value = '1-19%'
$ws.Range('AG47').Value = '1-19%'
$ws.Range('AJ47').Value = '1-19%'
$ws.Range('AE49').Value = '1-19%'
$ws.Range('AF49').Value = '1-19%'
$ws.Range('AG49').Value = '1-19%'
$ws.Range('AI49').Value = '1-19%'
$ws.Range('AJ49').Value = '1-19%'
$ws.Range('AF50').Value = '1-19%'
$ws.Range('AG51').Value = '1-19%'
$ws.Range('AH53').Value = '1-19%'
$ws.Range('AI53').Value = '1-19%'
$ws.Range('AE56').Value = '1-19%'
$ws.Range('AH56').Value = '1-19%'
$ws.Range('AF63').Value = '1-19%'
$ws.Range('AH63').Value = '1-19%'
$ws.Range('AE65').Value = '1-19%'
$ws.Range('AF65').Value = '1-19%'
$ws.Range('AG65').Value = '1-19%'
$ws.Range('AI65').Value = '1-19%'
$ws.Range('AJ65').Value = '1-19%'
$ws.Range('AH70').Value = '1-19%'
$ws.Range('AI70').Value = '1-19%'
$ws.Range('AE71').Value = '1-19%'
$ws.Range('AG71').Value = '1-19%'
$ws.Range('AH71').Value = '1-19%'
$ws.Range('AI71').Value = '1-19%'
$ws.Range('AJ71').Value = '1-19%'
$ws.Range('AE72').Value = '1-19%'
$ws.Range('AF72').Value = '1-19%'
$ws.Range('AG75').Value = '1-19%'
$ws.Range('AE76').Value = '1-19%'
$ws.Range('AF76').Value = '1-19%'
$ws.Range('AG76').Value = '1-19%'
$ws.Range('AH76').Value = '1-19%'
$ws.Range('AJ76').Value = '1-19%'
$ws.Range('AE77').Value = '1-19%'
$ws.Range('AI77').Value = '1-19%'
$ws.Range('AF78').Value = '1-19%'
$ws.Range('AH80').Value = '1-19%'
$ws.Range('AG82').Value = '1-19%'
$ws.Range('AI83').Value = '1-19%'
$ws.Range('AG84').Value = '1-19%'
$ws.Range('AF85').Value = '1-19%'
$ws.Range('AG85').Value = '1-19%'
$ws.Range('AE86').Value = '1-19%'
$ws.Range('AG87').Value = '1-19%'
$ws.Range('AJ87').Value = '1-19%'
$ws.Range('AG89').Value = '1-19%'
$ws.Range('AJ89').Value = '1-19%'
$ws.Range('AG90').Value = '1-19%'
$ws.Range('AE93').Value = '1-19%'
$ws.Range('AF93').Value = '1-19%'
$ws.Range('AG93').Value = '1-19%'
$ws.Range('AH93').Value = '1-19%'
$ws.Range('AI93').Value = '1-19%'
$ws.Range('AF94').Value = '1-19%'
$ws.Range('AE95').Value = '1-19%'
$ws.Range('AG95').Value = '1-19%'
$ws.Range('AH95').Value = '1-19%'
$ws.Range('AJ95').Value = '1-19%'
$ws.Range('AJ96').Value = '1-19%'
$ws.Range('AE97').Value = '1-19%'
$ws.Range('AH97').Value = '1-19%'
$ws.Range('AI97').Value = '1-19%'
$ws.Range('AJ97').Value = '1-19%'
$ws.Range('AG98').Value = '1-19%'
$ws.Range('AE99').Value = '1-19%'
$ws.Range('AI99').Value = '1-19%'
$ws.Range('AE100').Value = '1-19%'
$ws.Range('AF100').Value = '1-19%'
$ws.Range('AI100').Value = '1-19%'
$ws.Range('AH102').Value = '1-19%'
$ws.Range('AJ102').Value = '1-19%'
$ws.Range('AG103').Value = '1-19%'
$ws.Range('AH103').Value = '1-19%'
$ws.Range('AG104').Value = '1-19%'
$ws.Range('AE105').Value = '1-19%'
$ws.Range('AG105').Value = '1-19%'

# --- Clear stray numeric values from column AW for rows 16-70 ---
$ws.Range('AW16').ClearContents()
$ws.Range('AW17').ClearContents()
$ws.Range('AW18').ClearContents()
$ws.Range('AW19').ClearContents()
$ws.Range('AW20').ClearContents()
$ws.Range('AW21').ClearContents()
$ws.Range('AW23').ClearContents()
$ws.Range('AW24').ClearContents()
$ws.Range('AW26').ClearContents()
$ws.Range('AW28').ClearContents()
$ws.Range('AW29').ClearContents()
$ws.Range('AW30').ClearContents()
$ws.Range('AW31').ClearContents()
$ws.Range('AW32').ClearContents()
$ws.Range('AW33').ClearContents()
$ws.Range('AW34').ClearContents()
$ws.Range('AW35').ClearContents()
$ws.Range('AW36').ClearContents()
$ws.Range('AW37').ClearContents()
$ws.Range('AW39').ClearContents()
$ws.Range('AW40').ClearContents()
$ws.Range('AW41').ClearContents()
$ws.Range('AW42').ClearContents()
$ws.Range('AW43').ClearContents()
$ws.Range('AW44').ClearContents()
$ws.Range('AW45').ClearContents()
$ws.Range('AW46').ClearContents()
$ws.Range('AW47').ClearContents()
$ws.Range('AW48').ClearContents()
$ws.Range('AW50').ClearContents()
$ws.Range('AW51').ClearContents()
$ws.Range('AW52').ClearContents()
$ws.Range('AW53').ClearContents()
$ws.Range('AW54').ClearContents()
$ws.Range('AW55').ClearContents()
$ws.Range('AW56').ClearContents()
$ws.Range('AW57').ClearContents()
$ws.Range('AW59').ClearContents()
$ws.Range('AW60').ClearContents()
$ws.Range('AW62').ClearContents()
$ws.Range('AW63').ClearContents()
$ws.Range('AW64').ClearContents()
$ws.Range('AW65').ClearContents()
$ws.Range('AW66').ClearContents()
$ws.Range('AW67').ClearContents()
$ws.Range('AW68').ClearContents()
$ws.Range('AW69').ClearContents()
$ws.Range('AW70').ClearContents()

# --- View changes: zoom + selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 170
$ws.Range("AE105").Select()

# --- Column width changes for AE, AF, AG ---
$ws.Columns.Item(31).ColumnWidth = 14.0
$ws.Columns.Item(32).ColumnWidth = 16.0
$ws.Columns.Item(33).ColumnWidth = 13.333333333
